$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: force a run boundary at both ends of the given [start,end) range by
# round-tripping the Bold property (off/on or on/off, whichever is not the
# current value, then back). Any write to a character-formatting property
# causes the engine to re-normalise the run structure of the touched
# paragraph(s), and the touched sub-range always ends up as its own run once
# the write lands - exactly the "split" effect needed to reproduce the
# manually added <w:proofErr/> boundaries, since w:proofErr elements are not
# reachable through the exposed object model.
# ---------------------------------------------------------------------------
function SplitAt($posStart, $posEnd) {
    $rng = $d.Range($posStart, $posEnd)
    $cur = $rng.Bold
    $rng.Bold = 1 - $cur
    $rng.Bold = $cur
}

# ===========================================================================
# Edit 1 - "...downward facing ports(DFP)." gains a gramStart/gramEnd pair
# around "ports(" (no visible text change, just a run split at the word
# boundaries the grammar checker flagged).
# ===========================================================================
$r1 = $d.Content
$r1.Find.Execute("ports(DFP).") | Out-Null
$ports_start = $r1.Start
$ports_end = $ports_start + 6      # "ports(" is 6 characters
$dfp_end = $r1.End                 # end of "DFP)."

SplitAt $ports_start $ports_end
SplitAt $ports_end $dfp_end

# ===========================================================================
# Edit 2 - "... or  PI7C9X440SL" gains gramStart/gramEnd around "or  PI".
# Splits needed: " " | "or  " | "PI" | "7C9X440SL" (text itself is unchanged)
# ===========================================================================
$r2 = $d.Content
$r2.Find.Execute(" or  PI7C9X440SL") | Out-Null
$seg_start = $r2.Start
$split2a = $seg_start + 1
$split2b = $seg_start + 7

SplitAt $split2a $split2b

# ===========================================================================
# Edit 3 - "(i.e., MCDP2900) and a HDMI connector." run layout is rebuilt:
#   before: " " | "MCDP2900" | ")" | " " | "and a " | "HDMI" | " connector."
#   after : " MCDP2900)" | " " | "and a " | "HDMI" | " connector."
# Replacing the " MCDP2900)" span with itself merges the first three runs
# (the engine coalesces same-format adjacent runs on any mutating Find), then
# we re-split the boundaries that must stay distinct.
# ===========================================================================
$r3 = $d.Content
$r3.Find.Execute(" MCDP2900)", $true, $false, $false, $false, $false, $true, 1, $false, " MCDP2900)", 2) | Out-Null

$r3b = $d.Content
$r3b.Find.Execute(" MCDP2900) and a HDMI connector.") | Out-Null
$mcdp_start = $r3b.Start
$mcdp_close_end = $mcdp_start + 10     # end of " MCDP2900)"
$space_end = $mcdp_close_end + 1       # end of the lone " "
$anda_end = $space_end + 6             # end of "and a "
$hdmi_end = $anda_end + 4              # end of "HDMI"

SplitAt $mcdp_close_end $space_end
SplitAt $space_end $anda_end
SplitAt $anda_end $hdmi_end

# ===========================================================================
# Edit 4 - "The HD3SS460 will..." : "HD3SS460" + " " merge into one run
# "HD3SS460 ", while "The " and "will as well..." remain their own runs.
# ===========================================================================
$r4 = $d.Content
$r4.Find.Execute("The HD3SS460 will", $true, $false, $false, $false, $false, $true, 1, $false, "The HD3SS460 will", 2) | Out-Null

$r4b = $d.Content
$r4b.Find.Execute("The HD3SS460 will") | Out-Null
$the_start = $r4b.Start
$the_end = $the_start + 4          # end of "The "
$hd3_end = $the_end + 9            # end of "HD3SS460 "

SplitAt $the_start $the_end
SplitAt $the_end $hd3_end

# ===========================================================================
# Edit 5 - "TPS65982 from..." becomes "TPS65982" + new bold run
# "(we will the certified PD controller such as TPS65987) " + "from...".
# ===========================================================================
$r5 = $d.Content
$r5.Find.Execute("TPS65982 from") | Out-Null
$tps_start = $r5.Start
$space_start = $tps_start + 8
$space_end2 = $tps_start + 9

$insertedText = "(we will the certified PD controller such as TPS65987) "
$spaceRng = $d.Range($space_start, $space_end2)
$spaceRng.Text = $insertedText

$newRunEnd = $space_start + $insertedText.Length
SplitAt $space_start $newRunEnd
